# Applies scheduled market-price refresh updates to the Leve profit tables.
# For each affected row, H/I/J/K/L (price inputs) and M/N (profit outputs)
# are updated to the newly-scraped values. A few rows lose their N (HQ profit)
# cell entirely because the HQ price input column collapses to 0.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H81").Value = 52500
$ws.Range("J81").Value = 52500
$ws.Range("L81").Value = 52500
$ws.Range("N81").Value = -54496
$ws.Range("H84").Value = 52500
$ws.Range("J84").Value = 52500
$ws.Range("L84").Value = 157500
$ws.Range("N84").Value = -167484
$ws.Range("H129").Value = 366246.25
$ws.Range("J129").Value = 3648.4443
$ws.Range("L129").Value = 10945.3329
$ws.Range("N129").Value = -20945.3329
$ws.Range("H137").Value = 2793.5505
$ws.Range("I137").Value = 1031.04
$ws.Range("J137").Value = 3482.0312
$ws.Range("K137").Value = 3093.12
$ws.Range("L137").Value = 10446.0936
$ws.Range("M137").Value = -543.1199999999999
$ws.Range("N137").Value = -15546.0936

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2808.9
$ws.Range("I74").Value = 1600
$ws.Range("J74").Value = 3614.8333
$ws.Range("K74").Value = 1600
$ws.Range("L74").Value = 3614.8333
$ws.Range("M74").Value = -726
$ws.Range("N74").Value = -5362.8333
$ws.Range("H77").Value = 2808.9
$ws.Range("I77").Value = 1600
$ws.Range("J77").Value = 3614.8333
$ws.Range("K77").Value = 8000
$ws.Range("L77").Value = 18074.1665
$ws.Range("M77").Value = -3632
$ws.Range("N77").Value = -26810.1665
$ws.Range("H80").Value = 54674.4
$ws.Range("J80").Value = 54674.4
$ws.Range("L80").Value = 54674.4
$ws.Range("N80").Value = -56670.4
$ws.Range("H83").Value = 54674.4
$ws.Range("J83").Value = 54674.4
$ws.Range("L83").Value = 164023.2
$ws.Range("N83").Value = -174007.2
$ws.Range("H102").Value = 34868.168
$ws.Range("I102").Value = 1747.25
$ws.Range("J102").Value = 101110
$ws.Range("K102").Value = 1747.25
$ws.Range("L102").Value = 101110
$ws.Range("M102").Value = -125.25
$ws.Range("N102").Value = -104354

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H92").Value = 39996
$ws.Range("J92").Value = 39996
$ws.Range("L92").Value = 39996
$ws.Range("N92").Value = -44988
$ws.Range("H99").Value = 1701.238
$ws.Range("I99").Value = 1650.8334
$ws.Range("J99").Value = 2003.6666
$ws.Range("K99").Value = 1650.8334
$ws.Range("L99").Value = 2003.6666
$ws.Range("M99").Value = -152.8334
$ws.Range("N99").Value = -4999.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6945.878
$ws.Range("I31").Value = 4363.1665
$ws.Range("J31").Value = 7388.6284
$ws.Range("K31").Value = 4363.1665
$ws.Range("L31").Value = 7388.6284
$ws.Range("M31").Value = -4068.1665
$ws.Range("N31").Value = -7978.6284
$ws.Range("H34").Value = 6945.878
$ws.Range("I34").Value = 4363.1665
$ws.Range("J34").Value = 7388.6284
$ws.Range("K34").Value = 4363.1665
$ws.Range("L34").Value = 7388.6284
$ws.Range("M34").Value = -4161.1665
$ws.Range("N34").Value = -7792.6284
$ws.Range("H81").Value = 51312
$ws.Range("J81").Value = 51312
$ws.Range("L81").Value = 51312
$ws.Range("N81").Value = -53308
$ws.Range("H84").Value = 51312
$ws.Range("J84").Value = 51312
$ws.Range("L84").Value = 153936
$ws.Range("N84").Value = -163920
$ws.Range("H88").Value = 43267
$ws.Range("J88").Value = 43267
$ws.Range("L88").Value = 43267
$ws.Range("N88").Value = -44079
$ws.Range("H91").Value = 43267
$ws.Range("J91").Value = 43267
$ws.Range("L91").Value = 43267
$ws.Range("N91").Value = -46075
$ws.Range("H124").Value = 43764
$ws.Range("J124").Value = 43764
$ws.Range("L124").Value = 43764
$ws.Range("N124").Value = -48674
$ws.Range("H125").Value = 41079.5
$ws.Range("J125").Value = 41079.5
$ws.Range("L125").Value = 41079.5
$ws.Range("N125").Value = -45999.5
$ws.Range("H131").Value = 35663
$ws.Range("J131").Value = 35663
$ws.Range("L131").Value = 35663
$ws.Range("N131").Value = -45743
$ws.Range("H138").Value = 46800
$ws.Range("J138").Value = 46800
$ws.Range("L138").Value = 46800
$ws.Range("N138").Value = -57080

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H118").Value = 32458.8
$ws.Range("J118").Value = 32458.8
$ws.Range("L118").Value = 32458.8
$ws.Range("N118").Value = -35772.8
$ws.Range("H120").Value = 39309
$ws.Range("J120").Value = 39309
$ws.Range("L120").Value = 39309
$ws.Range("N120").Value = -48985
$ws.Range("H125").Value = 44326
$ws.Range("J125").Value = 44326
$ws.Range("L125").Value = 44326
$ws.Range("N125").Value = -49246
$ws.Range("H131").Value = 37324
$ws.Range("J131").Value = 37324
$ws.Range("L131").Value = 37324
$ws.Range("N131").Value = -47404

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()  # HQ input collapsed to 0; drop stale HQ-profit cell entirely
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()  # HQ input collapsed to 0; drop stale HQ-profit cell entirely
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()  # HQ input collapsed to 0; drop stale HQ-profit cell entirely
$ws.Range("H88").Value = 31511
$ws.Range("J88").Value = 44181
$ws.Range("L88").Value = 44181
$ws.Range("N88").Value = -45037
$ws.Range("H91").Value = 31511
$ws.Range("J91").Value = 44181
$ws.Range("L91").Value = 44181
$ws.Range("N91").Value = -47145
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()  # HQ input collapsed to 0; drop stale HQ-profit cell entirely
$ws.Range("H96").Value = 31299.25
$ws.Range("J96").Value = 31299.25
$ws.Range("L96").Value = 31299.25
$ws.Range("N96").Value = -36791.25
$ws.Range("H99").Value = 27166.666
$ws.Range("J99").Value = 29000
$ws.Range("L99").Value = 29000
$ws.Range("N99").Value = -34990
$ws.Range("H102").Value = 48553
$ws.Range("J102").Value = 48553
$ws.Range("L102").Value = 48553
$ws.Range("N102").Value = -55043
$ws.Range("H109").Value = 28346.334
$ws.Range("J109").Value = 28346.334
$ws.Range("L109").Value = 28346.334
$ws.Range("N109").Value = -31120.334
$ws.Range("H117").Value = 40384
$ws.Range("J117").Value = 40384
$ws.Range("L117").Value = 40384
$ws.Range("N117").Value = -49562
$ws.Range("H123").Value = 36429
$ws.Range("J123").Value = 36429
$ws.Range("L123").Value = 36429
$ws.Range("N123").Value = -46229
$ws.Range("H131").Value = 39996
$ws.Range("J131").Value = 39996
$ws.Range("L131").Value = 39996
$ws.Range("N131").Value = -50076
$ws.Range("H136").Value = 2969.8
$ws.Range("I136").Value = 2558.7
$ws.Range("J136").Value = 3792
$ws.Range("K136").Value = 7676.099999999999
$ws.Range("L136").Value = 11376
$ws.Range("M136").Value = -5126.099999999999
$ws.Range("N136").Value = -16476

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H27").Value = 31088.285
$ws.Range("J27").Value = 31088.285
$ws.Range("L27").Value = 31088.285
$ws.Range("N27").Value = -31226.285
$ws.Range("H93").Value = 38971.43
$ws.Range("J93").Value = 38971.43
$ws.Range("L93").Value = 38971.43
$ws.Range("N93").Value = -43963.43
$ws.Range("H102").Value = 31265.4
$ws.Range("J102").Value = 31265.4
$ws.Range("L102").Value = 31265.4
$ws.Range("N102").Value = -37755.4
$ws.Range("H109").Value = 33412
$ws.Range("J109").Value = 33412
$ws.Range("L109").Value = 33412
$ws.Range("N109").Value = -36186
$ws.Range("H115").Value = 37369
$ws.Range("J115").Value = 37369
$ws.Range("L115").Value = 37369
$ws.Range("N115").Value = -40503
$ws.Range("H118").Value = 33720.668
$ws.Range("J118").Value = 33720.668
$ws.Range("L118").Value = 33720.668
$ws.Range("N118").Value = -37034.668
$ws.Range("H127").Value = 43425
$ws.Range("J127").Value = 43425
$ws.Range("L127").Value = 43425
$ws.Range("N127").Value = -53345
$ws.Range("H129").Value = 31476
$ws.Range("J129").Value = 31476
$ws.Range("L129").Value = 31476
$ws.Range("N129").Value = -41476
$ws.Range("H136").Value = 19724.387
$ws.Range("I136").Value = 60240.59
$ws.Range("J136").Value = 2505
$ws.Range("K136").Value = 180721.77
$ws.Range("L136").Value = 7515
$ws.Range("M136").Value = -178171.77
$ws.Range("N136").Value = -12615
